$wb = $excel.ActiveWorkbook

# Sheet ALC row 33
$ws = $wb.Worksheets.Item("ALC")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 739
$arr[0,1] = 739
$arr[0,2] = 0
$arr[0,3] = 739
$arr[0,4] = 0
$arr[0,5] = -510
$arr[0,6] = $null
$ws.Range("H33:N33").Value = $arr

# Sheet ALC row 41
$ws = $wb.Worksheets.Item("ALC")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 420.6
$arr[0,1] = 429.2143
$arr[0,2] = 300
$arr[0,3] = 429.2143
$arr[0,4] = 300
$arr[0,5] = 10.78570000000002
$arr[0,6] = -1180
$ws.Range("H41:N41").Value = $arr

# Sheet ALC row 43
$ws = $wb.Worksheets.Item("ALC")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2749.625
$arr[0,1] = 2916
$arr[0,2] = 2250.5
$arr[0,3] = 2916
$arr[0,4] = 2250.5
$arr[0,5] = -2847
$arr[0,6] = -2388.5
$ws.Range("H43:N43").Value = $arr

# Sheet ALC row 76
$ws = $wb.Worksheets.Item("ALC")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2983
$arr[0,1] = 2974.5
$arr[0,2] = 3000
$arr[0,3] = 2974.5
$arr[0,4] = 3000
$arr[0,5] = -2659.5
$arr[0,6] = -3630
$ws.Range("H76:N76").Value = $arr

# Sheet ALC row 79
$ws = $wb.Worksheets.Item("ALC")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2983
$arr[0,1] = 2974.5
$arr[0,2] = 3000
$arr[0,3] = 2974.5
$arr[0,4] = 3000
$arr[0,5] = -1882.5
$arr[0,6] = -5184
$ws.Range("H79:N79").Value = $arr

# Sheet ALC row 111
$ws = $wb.Worksheets.Item("ALC")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 3396
$arr[0,1] = 3495
$arr[0,2] = 3000
$arr[0,3] = 10485
$arr[0,4] = 9000
$arr[0,5] = -7418
$arr[0,6] = -15134
$ws.Range("H111:N111").Value = $arr

# Sheet ALC row 129
$ws = $wb.Worksheets.Item("ALC")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2406.4285
$arr[0,1] = 1254.8334
$arr[0,2] = 3270.125
$arr[0,3] = 3764.5002
$arr[0,4] = 9810.375
$arr[0,5] = 1235.4998
$arr[0,6] = -19810.375
$ws.Range("H129:N129").Value = $arr

# Sheet ALC row 135
$ws = $wb.Worksheets.Item("ALC")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1809.6
$arr[0,1] = 1677.3334
$arr[0,2] = 3000
$arr[0,3] = 15096.0006
$arr[0,4] = 27000
$arr[0,5] = -12561.0006
$arr[0,6] = -32070
$ws.Range("H135:N135").Value = $arr

# Sheet ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2719.8
$arr[0,1] = 2649.75
$arr[0,2] = 3000
$arr[0,3] = 7949.25
$arr[0,4] = 9000
$arr[0,5] = -5399.25
$arr[0,6] = -14100
$ws.Range("H137:N137").Value = $arr

# Sheet ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1807.1333
$arr[0,1] = 1807.1333
$arr[0,2] = 0
$arr[0,3] = 1807.1333
$arr[0,4] = 0
$arr[0,5] = -1520.1333
$arr[0,6] = $null
$ws.Range("H32:N32").Value = $arr

# Sheet ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2869
$arr[0,1] = 2386.5
$arr[0,2] = 5764
$arr[0,3] = 2386.5
$arr[0,4] = 5764
$arr[0,5] = -2009.5
$arr[0,6] = -6518
$ws.Range("H45:N45").Value = $arr

# Sheet ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1154.1428
$arr[0,1] = 1089.0769
$arr[0,2] = 2000
$arr[0,3] = 3267.2307
$arr[0,4] = 6000
$arr[0,5] = -817.2307000000001
$arr[0,6] = -10900
$ws.Range("H122:N122").Value = $arr

# Sheet ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1676.1666
$arr[0,1] = 1601.2727
$arr[0,2] = 2500
$arr[0,3] = 4803.8181
$arr[0,4] = 7500
$arr[0,5] = -2273.8181
$arr[0,6] = -12560
$ws.Range("H132:N132").Value = $arr

# Sheet BSM row 29
$ws = $wb.Worksheets.Item("BSM")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 6666.6665
$arr[0,1] = 4000
$arr[0,2] = 8000
$arr[0,3] = 4000
$arr[0,4] = 8000
$arr[0,5] = -3711
$arr[0,6] = -8578
$ws.Range("H29:N29").Value = $arr

# Sheet BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 4000
$arr[0,1] = 0
$arr[0,2] = 4000
$arr[0,3] = 0
$arr[0,4] = 4000
$arr[0,5] = $null
$arr[0,6] = -6996
$ws.Range("H99:N99").Value = $arr

# Sheet CRP row 16
$ws = $wb.Worksheets.Item("CRP")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1716.3334
$arr[0,1] = 1716.3334
$arr[0,2] = 0
$arr[0,3] = 1716.3334
$arr[0,4] = 0
$arr[0,5] = -1429.3334
$arr[0,6] = $null
$ws.Range("H16:N16").Value = $arr

# Sheet CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 6166.1665
$arr[0,1] = 5249.25
$arr[0,2] = 8000
$arr[0,3] = 5249.25
$arr[0,4] = 8000
$arr[0,5] = -5046.25
$arr[0,6] = -8406
$ws.Range("H58:N58").Value = $arr

# Sheet CRP row 62
$ws = $wb.Worksheets.Item("CRP")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 3966.3333
$arr[0,1] = 1900
$arr[0,2] = 4999.5
$arr[0,3] = 1900
$arr[0,4] = 4999.5
$arr[0,5] = -1276
$arr[0,6] = -6247.5
$ws.Range("H62:N62").Value = $arr

# Sheet CRP row 65
$ws = $wb.Worksheets.Item("CRP")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 3966.3333
$arr[0,1] = 1900
$arr[0,2] = 4999.5
$arr[0,3] = 9500
$arr[0,4] = 24997.5
$arr[0,5] = -6380
$arr[0,6] = -31237.5
$ws.Range("H65:N65").Value = $arr

# Sheet CRP row 86
$ws = $wb.Worksheets.Item("CRP")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 3381.6667
$arr[0,1] = 3379.2
$arr[0,2] = 3394
$arr[0,3] = 3379.2
$arr[0,4] = 3394
$arr[0,5] = -2256.2
$arr[0,6] = -5640
$ws.Range("H86:N86").Value = $arr

# Sheet CRP row 89
$ws = $wb.Worksheets.Item("CRP")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 3381.6667
$arr[0,1] = 3379.2
$arr[0,2] = 3394
$arr[0,3] = 16896
$arr[0,4] = 16970
$arr[0,5] = -11280
$arr[0,6] = -28202
$ws.Range("H89:N89").Value = $arr

# Sheet CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 3828.7144
$arr[0,1] = 3760.2
$arr[0,2] = 4000
$arr[0,3] = 3760.2
$arr[0,4] = 4000
$arr[0,5] = -2262.2
$arr[0,6] = -6996
$ws.Range("H99:N99").Value = $arr

# Sheet CRP row 113
$ws = $wb.Worksheets.Item("CRP")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1716.3334
$arr[0,1] = 1716.3334
$arr[0,2] = 0
$arr[0,3] = 1716.3334
$arr[0,4] = 0
$arr[0,5] = 453.6666
$arr[0,6] = $null
$ws.Range("H113:N113").Value = $arr

# Sheet CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 3828.7144
$arr[0,1] = 3760.2
$arr[0,2] = 4000
$arr[0,3] = 11280.6
$arr[0,4] = 12000
$arr[0,5] = -8810.599999999999
$arr[0,6] = -16940
$ws.Range("H126:N126").Value = $arr

# Sheet CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 6166.1665
$arr[0,1] = 5249.25
$arr[0,2] = 8000
$arr[0,3] = 15747.75
$arr[0,4] = 24000
$arr[0,5] = -13197.75
$arr[0,6] = -29100
$ws.Range("H136:N136").Value = $arr

# Sheet CUL row 4
$ws = $wb.Worksheets.Item("CUL")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2850644.5
$arr[0,1] = 1062639
$arr[0,2] = 10002666
$arr[0,3] = 3187917
$arr[0,4] = 30007998
$arr[0,5] = -3187805
$arr[0,6] = -30008222
$ws.Range("H4:N4").Value = $arr

# Sheet CUL row 9
$ws = $wb.Worksheets.Item("CUL")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2280.4
$arr[0,1] = 1600.5
$arr[0,2] = 5000
$arr[0,3] = 4801.5
$arr[0,4] = 15000
$arr[0,5] = -4577.5
$arr[0,6] = -15448
$ws.Range("H9:N9").Value = $arr

# Sheet CUL row 121
$ws = $wb.Worksheets.Item("CUL")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 59413652
$arr[0,1] = 1061.2858
$arr[0,2] = 101002460
$arr[0,3] = 3183.8574
$arr[0,4] = 303007380
$arr[0,5] = -1873.8574
$arr[0,6] = -303010000
$ws.Range("H121:N121").Value = $arr

# Sheet GSM row 97
$ws = $wb.Worksheets.Item("GSM")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 868.1177
$arr[0,1] = 857.8333
$arr[0,2] = 892.8
$arr[0,3] = 857.8333
$arr[0,4] = 892.8
$arr[0,5] = -361.8333
$arr[0,6] = -1884.8
$ws.Range("H97:N97").Value = $arr

# Sheet GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 3045
$arr[0,1] = 3045
$arr[0,2] = 0
$arr[0,3] = 3045
$arr[0,4] = 0
$arr[0,5] = -1423
$arr[0,6] = $null
$ws.Range("H102:N102").Value = $arr

# Sheet GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 4612.6
$arr[0,1] = 4513.75
$arr[0,2] = 5008
$arr[0,3] = 13541.25
$arr[0,4] = 15024
$arr[0,5] = -11091.25
$arr[0,6] = -19924
$ws.Range("H122:N122").Value = $arr

# Sheet GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 1709.9048
$arr[0,1] = 1709.9048
$arr[0,2] = 0
$arr[0,3] = 5129.7144
$arr[0,4] = 0
$arr[0,5] = -2599.7144
$arr[0,6] = $null
$ws.Range("H132:N132").Value = $arr

# Sheet LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 3389
$arr[0,1] = 4185.3335
$arr[0,2] = 1000
$arr[0,3] = 4185.3335
$arr[0,4] = 1000
$arr[0,5] = -3890.3335
$arr[0,6] = -1590
$ws.Range("H22:N22").Value = $arr

# Sheet LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 3389
$arr[0,1] = 4185.3335
$arr[0,2] = 1000
$arr[0,3] = 4185.3335
$arr[0,4] = 1000
$arr[0,5] = -4078.3335
$arr[0,6] = -1214
$ws.Range("H27:N27").Value = $arr

# Sheet LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 4104.4
$arr[0,1] = 4223.857
$arr[0,2] = 3999.875
$arr[0,3] = 4223.857
$arr[0,4] = 3999.875
$arr[0,5] = -4035.857
$arr[0,6] = -4375.875
$ws.Range("H46:N46").Value = $arr

# Sheet LTW row 68
$ws = $wb.Worksheets.Item("LTW")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 4804.25
$arr[0,1] = 5739
$arr[0,2] = 2000
$arr[0,3] = 5739
$arr[0,4] = 2000
$arr[0,5] = -4990
$arr[0,6] = -3498
$ws.Range("H68:N68").Value = $arr

# Sheet LTW row 71
$ws = $wb.Worksheets.Item("LTW")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 4804.25
$arr[0,1] = 5739
$arr[0,2] = 2000
$arr[0,3] = 28695
$arr[0,4] = 10000
$arr[0,5] = -24951
$arr[0,6] = -17488
$ws.Range("H71:N71").Value = $arr

# Sheet LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2574.2188
$arr[0,1] = 2637.6365
$arr[0,2] = 2541
$arr[0,3] = 7912.9095
$arr[0,4] = 7623
$arr[0,5] = -5462.9095
$arr[0,6] = -12523
$ws.Range("H122:N122").Value = $arr

# Sheet LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 3334.05
$arr[0,1] = 3461.3125
$arr[0,2] = 2825
$arr[0,3] = 10383.9375
$arr[0,4] = 8475
$arr[0,5] = -7853.9375
$arr[0,6] = -13535
$ws.Range("H132:N132").Value = $arr

# Sheet LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 14583.167
$arr[0,1] = 12000
$arr[0,2] = 19749.5
$arr[0,3] = 36000
$arr[0,4] = 59248.5
$arr[0,5] = -33450
$arr[0,6] = -64348.5
$ws.Range("H136:N136").Value = $arr

# Sheet WVR row 15
$ws = $wb.Worksheets.Item("WVR")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 10000
$arr[0,1] = 10000
$arr[0,2] = 0
$arr[0,3] = 10000
$arr[0,4] = 0
$arr[0,5] = -9712
$arr[0,6] = $null
$ws.Range("H15:N15").Value = $arr

# Sheet WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 2567.2
$arr[0,1] = 2296.3333
$arr[0,2] = 5005
$arr[0,3] = 6888.999899999999
$arr[0,4] = 15015
$arr[0,5] = -4438.999899999999
$arr[0,6] = -19915
$ws.Range("H122:N122").Value = $arr

# Sheet WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$arr = New-Object 'object[,]' 1,7
$arr[0,0] = 10775
$arr[0,1] = 10775
$arr[0,2] = 0
$arr[0,3] = 32325
$arr[0,4] = 0
$arr[0,5] = -29775
$arr[0,6] = $null
$ws.Range("H136:N136").Value = $arr
